# "Version finale du PowerPoint"
#
# The "Sommaire" slide (slide 2) holds a SmartArt list diagram whose
# items II and IV had their labels swapped by mistake:
#   II. Présentation du SWOT        -> should be  II. Présentation du QUINTILIEN
#   IV. Présentation du QUINTILIEN  -> should be  IV. Présentation du SWOT
#
# Fix the two SmartArt node labels directly through the SmartArt object
# model so both the diagram data and its drawing cache stay in sync.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the shape that hosts the SmartArt graphic on the slide.
$smartArtShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasSmartArt) {
        $smartArtShape = $candidate
        break
    }
}

$nodes = $smartArtShape.SmartArt.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $text = $node.TextFrame2.TextRange.Text

    if ($text -eq "II. Présentation du SWOT") {
        $node.TextFrame2.TextRange.Text = "II. Présentation du QUINTILIEN"
    }
    elseif ($text -eq "IV. Présentation du QUINTILIEN") {
        $node.TextFrame2.TextRange.Text = "IV. Présentation du SWOT"
    }
}
